$wb = $excel.ActiveWorkbook

# --- weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.28927585322432
$ws.Range("C2").Value = 0.263546347264198
$ws.Range("B3").Value = -0.107253607990594
$ws.Range("C3").Value = 0.169984735668374

# --- lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 1.52982703612423
$ws.Range("C2").Value = 0.343917875306103
$ws.Range("B3").Value = -0.830115526445352
$ws.Range("C3").Value = 0.186708829854108

# --- llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.83920081528535
$ws.Range("C2").Value = 0.283458625754857
$ws.Range("B3").Value = 1.43323566724509
$ws.Range("C3").Value = 0.380173870798206

# --- gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.15271197631863
$ws.Range("C2").Value = 0.295240845282059
$ws.Range("B3").Value = -0.0312207305439134
$ws.Range("C3").Value = 0.0287361000618847

# --- exp: unchanged ---

# --- weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0694566771563014
$ws.Range("B2").Value = -0.0143544642403251
$ws.Range("A3").Value = -0.0143544642403251
$ws.Range("B3").Value = 0.0288948103602471

# --- lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.118279504955064
$ws.Range("B2").Value = -0.0477083810907798
$ws.Range("A3").Value = -0.0477083810907798
$ws.Range("B3").Value = 0.0348601871454903

# --- llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0803487925148324
$ws.Range("B2").Value = 0.0219698548481498
$ws.Range("A3").Value = 0.0219698548481498
$ws.Range("B3").Value = 0.144532172037691

# --- gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0871671567228649
$ws.Range("B2").Value = -0.00352109315215479
$ws.Range("A3").Value = -0.00352109315215479
$ws.Range("B3").Value = 0.000825763446766651

# --- exp cov: unchanged ---
